$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 708.5714
$ws.Range("I33").Value = 491.75
$ws.Range("J33").Value = 997.6667
$ws.Range("K33").Value = 491.75
$ws.Range("L33").Value = 997.6667
$ws.Range("M33").Value = -262.75
$ws.Range("N33").Value = -1455.6667

$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1825

$ws.Range("H125").Value = 10007
$ws.Range("I125").Value = 5000
$ws.Range("K125").Value = 45000
$ws.Range("M125").Value = -42540

$ws.Range("H129").Value = 3336.077
$ws.Range("I129").Value = 1323.75
$ws.Range("K129").Value = 3971.25
$ws.Range("M129").Value = 1028.75

$ws.Range("H131").Value = 1071.625
$ws.Range("I131").Value = 796.1429000000001
$ws.Range("J131").Value = 3000
$ws.Range("K131").Value = 2388.4287
$ws.Range("L131").Value = 9000
$ws.Range("M131").Value = 2651.5713
$ws.Range("N131").Value = -19080

$ws.Range("H138").Value = 1669.3478
$ws.Range("I138").Value = 479
$ws.Range("K138").Value = 1437
$ws.Range("M138").Value = 3703

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2158.2666
$ws.Range("I2").Value = 2352.7273
$ws.Range("K2").Value = 2352.7273
$ws.Range("M2").Value = -2239.7273

$ws.Range("H44").Value = 33999.5
$ws.Range("I44").Value = 32999
$ws.Range("K44").Value = 32999
$ws.Range("M44").Value = -32511

$ws.Range("H45").Value = 1803
$ws.Range("I45").Value = 1803
$ws.Range("K45").Value = 1803
$ws.Range("M45").Value = -1426

$ws.Range("H55").Value = 23999.666
$ws.Range("I55").Value = 21999
$ws.Range("K55").Value = 21999
$ws.Range("M55").Value = -21684

$ws.Range("H80").Value = 35713.57
$ws.Range("J80").Value = 39999.168
$ws.Range("L80").Value = 39999.168
$ws.Range("N80").Value = -41995.168

$ws.Range("H83").Value = 35713.57
$ws.Range("J83").Value = 39999.168
$ws.Range("L83").Value = 119997.504
$ws.Range("N83").Value = -129981.504

$ws.Range("H116").Value = 2158.2666
$ws.Range("I116").Value = 2352.7273
$ws.Range("K116").Value = 2352.7273
$ws.Range("M116").Value = -58.72730000000001

$ws.Range("H122").Value = 2997.25
$ws.Range("I122").Value = 2996.3333
$ws.Range("K122").Value = 8988.999899999999
$ws.Range("M122").Value = -6538.999899999999

$ws.Range("H130").Value = 19999.5
$ws.Range("J130").Value = 19999.5
$ws.Range("L130").Value = 19999.5
$ws.Range("N130").Value = -30039.5

$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()

$ws.Range("H132").Value = 1394.5
$ws.Range("I132").Value = 1348.2916
$ws.Range("J132").Value = 1671.75
$ws.Range("K132").Value = 4044.8748
$ws.Range("L132").Value = 5015.25
$ws.Range("M132").Value = -1514.8748
$ws.Range("N132").Value = -10075.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2158.2666
$ws.Range("I3").Value = 2352.7273
$ws.Range("K3").Value = 2352.7273
$ws.Range("M3").Value = -2238.7273

$ws.Range("H20").Value = 1259.8462
$ws.Range("I20").Value = 661.8
$ws.Range("J20").Value = 3253.3333
$ws.Range("K20").Value = 661.8
$ws.Range("L20").Value = 3253.3333
$ws.Range("M20").Value = -414.8
$ws.Range("N20").Value = -3747.3333

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H99").Value = 7272.6665
$ws.Range("I99").Value = 7272.6665
$ws.Range("K99").Value = 7272.6665
$ws.Range("M99").Value = -5774.6665

$ws.Range("H126").Value = 7272.6665
$ws.Range("I126").Value = 7272.6665
$ws.Range("K126").Value = 21817.9995
$ws.Range("M126").Value = -19347.9995

$ws.Range("H132").Value = 1542.6
$ws.Range("I132").Value = 1695.6666
$ws.Range("K132").Value = 5086.9998
$ws.Range("M132").Value = -2556.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 115
$ws.Range("J52").Value = 115
$ws.Range("L52").Value = 345
$ws.Range("N52").Value = -877

$ws.Range("H56").Value = 10000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 10000
$ws.Range("M56").Value = -9470

$ws.Range("H81").Value = 11380.5
$ws.Range("J81").Value = 12720.571
$ws.Range("L81").Value = 38161.713
$ws.Range("N81").Value = -40407.713

$ws.Range("H84").Value = 11380.5
$ws.Range("J84").Value = 12720.571
$ws.Range("L84").Value = 114485.139
$ws.Range("N84").Value = -125717.139

$ws.Range("H102").Value = 7999
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws.Range("H104").Value = 26666.666
$ws.Range("I104").Value = 15500
$ws.Range("K104").Value = 46500
$ws.Range("M104").Value = -43879

$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()

$ws.Range("H129").Value = 500
$ws.Range("I129").Value = 500
$ws.Range("K129").Value = 1500
$ws.Range("M129").Value = 3500

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 4440
$ws.Range("I41").Value = 3300
$ws.Range("J41").Value = 9000
$ws.Range("K41").Value = 3300
$ws.Range("L41").Value = 9000
$ws.Range("M41").Value = -2945
$ws.Range("N41").Value = -9710

$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H122").Value = 31254918
$ws.Range("J122").Value = 6668
$ws.Range("L122").Value = 20004
$ws.Range("N122").Value = -24904

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5458.3335
$ws.Range("I7").Value = 4229
$ws.Range("J7").Value = 6995
$ws.Range("K7").Value = 4229
$ws.Range("L7").Value = 6995
$ws.Range("M7").Value = -4117
$ws.Range("N7").Value = -7219

$ws.Range("H30").Value = 837.25
$ws.Range("J30").Value = 1500
$ws.Range("L30").Value = 1500
$ws.Range("N30").Value = -1716

$ws.Range("H46").Value = 3624.8333
$ws.Range("J46").Value = 3833.111
$ws.Range("L46").Value = 3833.111
$ws.Range("N46").Value = -4209.111

$ws.Range("H61").Value = 5097.6665
$ws.Range("I61").Value = 5097.6665
$ws.Range("K61").Value = 5097.6665
$ws.Range("M61").Value = -4895.6665

$ws.Range("H113").Value = 5097.6665
$ws.Range("I113").Value = 5097.6665
$ws.Range("K113").Value = 5097.6665
$ws.Range("M113").Value = -2927.6665

$ws.Range("H126").Value = 5458.3335
$ws.Range("I126").Value = 4229
$ws.Range("J126").Value = 6995
$ws.Range("K126").Value = 12687
$ws.Range("L126").Value = 20985
$ws.Range("M126").Value = -10217
$ws.Range("N126").Value = -25925

$ws.Range("H131").Value = 16992
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws.Range("H132").Value = 5988.5386
$ws.Range("I132").Value = 5806.9546
$ws.Range("K132").Value = 17420.8638
$ws.Range("M132").Value = -14890.8638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 10207.333
$ws.Range("J20").Value = 10207.333
$ws.Range("L20").Value = 10207.333
$ws.Range("N20").Value = -10687.333

$ws.Range("H22").Value = 3200
$ws.Range("J22").Value = 3200
$ws.Range("L22").Value = 3200
$ws.Range("N22").Value = -3786

$ws.Range("H23").Value = 14485.571
$ws.Range("I23").Value = 20159.8
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 20159.8
$ws.Range("L23").Value = 300
$ws.Range("M23").Value = -19930.8
$ws.Range("N23").Value = -758

$ws.Range("H30").Value = 4009
$ws.Range("I30").Value = 4009
$ws.Range("K30").Value = 4009
$ws.Range("M30").Value = -3902

$ws.Range("H41").Value = 15956.7
$ws.Range("I41").Value = 16724.5
$ws.Range("J41").Value = 15444.833
$ws.Range("K41").Value = 16724.5
$ws.Range("L41").Value = 15444.833
$ws.Range("M41").Value = -16334.5
$ws.Range("N41").Value = -16224.833

$ws.Range("H126").Value = 3012.35
$ws.Range("I126").Value = 3028.1667
$ws.Range("J126").Value = 2988.625
$ws.Range("K126").Value = 9084.500100000001
$ws.Range("L126").Value = 8965.875
$ws.Range("M126").Value = -6614.500100000001
$ws.Range("N126").Value = -13905.875

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H131").Value = 30650
$ws.Range("I131").Value = 30650
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 30650
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -25610
$ws.Range("N131").ClearContents()

$ws.Range("H132").Value = 1213.4286
$ws.Range("I132").Value = 1199
$ws.Range("J132").Value = 1249.5
$ws.Range("K132").Value = 3597
$ws.Range("L132").Value = 3748.5
$ws.Range("M132").Value = -1067
$ws.Range("N132").Value = -8808.5
